# Automatische test-sync: 2025-06-19 17:53:30
#
# Adds the new incoming mail-log entry (row 27) to the "Logs" sheet and
# refreshes the "Dashboard" pivot-style summary so the "Klacht" category
# (now tied with "Bestelling" at 3) is reflected and re-ordered ahead of
# "Bestelling".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new row.
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 27
$logs.Cells.Item($newRow, 1).Value = "Klacht over levering"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$logs.Cells.Item($newRow, 4).Value = "Klacht"
$logs.Cells.Item($newRow, 6).Value = "2025-06-19 17:53:20"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Extend the existing conditional-formatting rule groups (Categorie / D,
# and Beantwoord / G) to cover the newly added row while keeping their
# rules, priorities and dxf formats untouched.
$logs.Range("D2:D26").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D27"))
$logs.Range("G2:G26").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G27"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: "Klacht" now also totals 3 (tied with
#    "Bestelling"), so it moves up to row 4 and "Bestelling" drops to
#    row 5, with its count updated to 3.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Klacht"
$dash.Cells.Item(5, 1).Value = "Bestelling"
$dash.Cells.Item(5, 2).Value = 3
